# Weekly fruit/vegetable price update: insert one new observation row
# at row 78 (pushing the existing rows 78-103 down to 79-104), matching
# the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 78; everything below shifts down
# by one (old row 78 -> 79, ..., old row 103 -> 104).
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly record.
$ws.Cells.Item(78, 1).Value = 4
$ws.Cells.Item(78, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(78, 3).Value = 'Los Lagos'
$ws.Cells.Item(78, 4).Value = 44841
$ws.Cells.Item(78, 5).Value = 10
$ws.Cells.Item(78, 6).Value = 100112031
$ws.Cells.Item(78, 7).Value = 'Poroto verde'
$ws.Cells.Item(78, 8).Value = 'Magnum'
$ws.Cells.Item(78, 9).Value = 'Primera'
$ws.Cells.Item(78, 10).Value = 35
$ws.Cells.Item(78, 11).Value = 35000
$ws.Cells.Item(78, 12).Value = 35000
$ws.Cells.Item(78, 13).Value = 35000
$ws.Cells.Item(78, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(78, 15).Value = 'Perú'
$ws.Cells.Item(78, 16).Value = 1400
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = 'Hortaliza'
